$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScoutingData")

$ws.Range("A9").Value = "f"
$ws.Range("B9").Value = "g"
$ws.Range("C9").Value = "qg"
$ws.Range("D9").Value = "fd"
$ws.Range("E9").Value = "fdsafdsa hi"
